# Generate Report for Handoff
#
# The e9de8414-4b79-4776-a2b2-0975ebd1684d file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff" - update the
# Overview roll-up plus the per-locale (zh-cn / de-de) detail sheets with
# the new status, the refreshed timestamps and the stale-handback error
# message, and widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$newOverviewDt  = "2016-08-22 12:50:06"
$zhHandoffDt    = "2016-08-22 12:49:56"
$errorDetail    = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ac5de5c69375ac589e7ed2d8f998a18e26550c20/e2e/e9de8414-4b79-4776-a2b2-0975ebd1684d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d659163232932c9fa71f859f16b3c2c749bd9c9a/e2e/e9de8414-4b79-4776-a2b2-0975ebd1684d.md."

# --- Overview sheet: row 3 is the e9de8414...md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $newOverviewDt

# --- zh-cn detail sheet: row 3 is the e9de8414...md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("H3").Value = $zhHandoffDt
$zhcn.Range("P3").Value = $errorDetail
# Widen the Error Detail column (P) so the long message is visible.
# ColumnWidth is character-width based and Excel stores it in the package
# with ~0.8333 char padding added, so back that padding out to land on an
# exact stored width of 40.
$zhcn.Range("P1").ColumnWidth = 39.166666666666664

# --- de-de detail sheet: row 3 is the e9de8414...md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("H3").Value = $newOverviewDt
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.166666666666664
